$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking snapshot refresh: prices + 1h volume deltas for all 50 rows,
# plus a ranking swap for two pairs of coins (Stellar/TheGraph, ARBITRUM/NEARProtocol).
# Force text storage (matches existing inlineStr cells) so values such as "1.00" or
# "0.0000218" are not silently reinterpreted as numbers by the Value setter.
$textRange = $ws.Range("B2:E51")
$textRange.NumberFormat = "@"

$ws.Range("D2").Value = "62.723.85"
$ws.Range("E2").Value = "  +1.09%  "
$ws.Range("D3").Value = "3.466.73"
$ws.Range("E3").Value = "  +1.12%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "414.09"
$ws.Range("E5").Value = "  +1.10%  "
$ws.Range("D6").Value = "130.30"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").Value = "0.626"
$ws.Range("E7").Value = "  -1.37%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.726"
$ws.Range("E9").Value = "  -1.98%  "
$ws.Range("D10").Value = "0.148"
$ws.Range("E10").Value = "  +4.37%  "
$ws.Range("D11").Value = "42.54"
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("D12").Value = "9.61"
$ws.Range("E12").Value = "  +3.90%  "
$ws.Range("D13").Value = "0.0000218"
$ws.Range("E13").Value = "  -2.98%  "
$ws.Range("D14").Value = "4.019.57"
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("E16").Value = "  -3.89%  "
$ws.Range("D17").Value = "3.480.94"
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("D18").Value = "12.67"
$ws.Range("E18").Value = "  +1.30%  "
$ws.Range("E19").Value = "  -1.77%  "
$ws.Range("D20").Value = "62.703.28"
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("D21").Value = "462.28"
$ws.Range("E21").Value = "  +1.70%  "
$ws.Range("D22").Value = "90.56"
$ws.Range("E22").Value = "  -1.01%  "
$ws.Range("E23").Value = "  +1.60%  "
$ws.Range("D24").Value = "13.29"
$ws.Range("E24").Value = "  +1.18%  "
$ws.Range("D25").Value = "10.75"
$ws.Range("E25").Value = "  +17.97%  "
$ws.Range("D26").Value = "3.32"
$ws.Range("E26").Value = "  +1.37%  "
$ws.Range("D27").Value = "33.33"
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("E28").Value = "  +0.36%  "
$ws.Range("D29").Value = "7.56"
$ws.Range("E29").Value = "  -1.21%  "
$ws.Range("E30").Value = "  -1.01%  "
$ws.Range("E32").Value = "  -2.66%  "
$ws.Range("E33").Value = "  -2.18%  "
$ws.Range("D34").Value = "40.72"
$ws.Range("E34").Value = "  -5.33%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "58.46"
$ws.Range("E36").Value = "  +7.50%  "
$ws.Range("E38").Value = "  +4.62%  "
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("D40").Value = "147.65"
$ws.Range("E40").Value = "  +3.61%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").Value = "0.134"
$ws.Range("E41").Value = "  -0.74%  "
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").Value = "0.320"
$ws.Range("E42").Value = "  +0.28%  "
$ws.Range("E43").Value = "  -1.50%  "
$ws.Range("E44").Value = "  +5.64%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "4.36"
$ws.Range("E45").Value = "  +1.85%  "
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").Value = "2.07"
$ws.Range("E46").Value = "  +3.48%  "
$ws.Range("D47").Value = "2.41"
$ws.Range("E47").Value = "  +13.86%  "
$ws.Range("D48").Value = "0.0₃0557"
$ws.Range("E48").Value = "  +30.41%  "
$ws.Range("E49").Value = "  -1.59%  "
$ws.Range("D50").Value = "22.20"
$ws.Range("E50").Value = "  -1.01%  "
$ws.Range("D51").Value = "0.141"
$ws.Range("E51").Value = "  +0.71%  "

$textRange.ClearFormats()
